$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The scraper now pulls two extra columns (height, weight) for every
# player row. They land right after the existing stat columns (A-D) and
# before the old "fantasy points" column, so the previous column-E
# ("fantasy points") values have to shift two columns over to G.
$fantasyPoints = @{
    2  = 1.5
    3  = 1
    4  = 0.7
    5  = 5.7
    6  = 1.8
    7  = 2.3
    8  = 0
    9  = 6.4
    10 = 14.3
    11 = 2
    12 = 4.4
}

foreach ($r in $fantasyPoints.Keys) {
    $ws.Cells.Item($r, 7).Value = $fantasyPoints[$r]
}

# Populate the new height (E) and weight (F) columns for every data row.
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 5).Value = 6.416666666666667
    $ws.Cells.Item($r, 6).Value = 255
}

# Headers: E1/F1 are the new "height"/"weight" columns, and G1 becomes the
# relocated "fantasy points" header.
$ws.Cells.Item(1, 5).Value = "height"
$ws.Cells.Item(1, 6).Value = "weight"
$ws.Cells.Item(1, 7).Value = "fantasy points"

# Match the header formatting (bold, centered, bordered) used by the rest
# of row 1.
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1:G1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
